$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Rent amount (was 300, now 2300)
$ws.Range("B2").Value = 2300

# Add two new expense rows: Travel and Groceries
$ws.Range("A3").Value = "Travel"
$ws.Range("B3").Value = 380
$ws.Range("C3").Value = 46060.291712962964

$ws.Range("A4").Value = "Groceries"
$ws.Range("B4").Value = 430
$ws.Range("C4").Value = 46051.291712962964

# Match the date formatting already used on C2
$ws.Range("C2").Copy()
$ws.Range("C3:C4").PasteSpecial(-4122)
